$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height adjustments (rows whose customHeight changes) ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30

# --- Clear cells whose old content is not part of the final layout ---
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()

# --- Cell content updates (final layout) ---
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'
$ws.Range("B2").Value = 'LOQ4078'
$ws.Range("C2").Value = 'LOQ4078'
$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Modelagem e Simulação de Processos'
$ws.Range("C3").Value = ' Modelagem e Simulação de Processos'
$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Process Modeling and Simulation'
$ws.Range("C4").Value = 'Process Modeling and Simulation'
$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '4'
$ws.Range("C5").Value = '4'
$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'
$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'
$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2018'
$ws.Range("C8").Value = '01/01/2018'
$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EQD-8,EQN-10'
$ws.Range("C9").Value = 'EQD-8,EQN-10'
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '5817066 - Félix Monteiro Pereira'
$ws.Range("C10").Value = '5817066 - Félix Monteiro Pereira'
$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'This course aims to introduce the students the modeling and simulation tools present in the work environment of the chemical engineer. These tools aids to design, operation and optimization of industrial processes.'
$ws.Range("C11").Value = 'This course aims to introduce the students the modeling and simulation tools present in the work environment of the chemical engineer. These tools aids to design, operation and optimization of industrial processes.'
$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Issues of interest and concern of Chemical Engineers will be addressed with an emphasis on advances and innovations of basic aspects of the Chemical Engineering program.'
$ws.Range("C14").Value = 'Issues of interest and concern of Chemical Engineers will be addressed with an emphasis on advances and innovations of basic aspects of the Chemical Engineering program.'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2018'
$ws.Range("C15").Value = '01/01/2018'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1. Introduction to process modeling and simulation. Problem solutions in chemical engineering with: systems of algebraic linear equations, systems of non-linear algebraic equations, systems of first-order differential equations (initial value problems), second order differential equations (boundary value problems), partial differential equations, process optimization, Non-linear regression and statistical analysis.2. Introduction to process simulators: logical operations, unit operations, chemical reaction units, process flowchart.'
$ws.Range("C16").Value = '1. Introduction to process modeling and simulation. Problem solutions in chemical engineering with: systems of algebraic linear equations, systems of non-linear algebraic equations, systems of first-order differential equations (initial value problems), second order differential equations (boundary value problems), partial differential equations, process optimization, Non-linear regression and statistical analysis.2. Introduction to process simulators: logical operations, unit operations, chemical reaction units, process flowchart.'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '5817066 - Félix Monteiro Pereira'
$ws.Range("C18").Value = '5817066 - Félix Monteiro Pereira'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'A avaliação do aprendizado será realizada por meio da resolução de problemas de engenharia química propostos aos alunos. A média final será calculada pela média aritmética entre a nota obtida na resolução de problemas relativos aos itens 1 e 2 do programa do curso.'
$ws.Range("C19").Value = 'A avaliação do aprendizado será realizada por meio da resolução de problemas de engenharia química propostos aos alunos. A média final será calculada pela média aritmética entre a nota obtida na resolução de problemas relativos aos itens 1 e 2 do programa do curso.'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'MF = (P1 + P2)/2 Onde: P1 é a nota obtida pela avaliação da resolução de problemas referentes ao item 1 do Programa do curso;P2 é a nota obtida pela avaliação da resolução de problemas referentes ao item 2 do Programa do curso;MF é a média final do período.'
$ws.Range("C20").Value = 'MF = (P1 + P2)/2 Onde: P1 é a nota obtida pela avaliação da resolução de problemas referentes ao item 1 do Programa do curso;P2 é a nota obtida pela avaliação da resolução de problemas referentes ao item 2 do Programa do curso;MF é a média final do período.'
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0.'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0.'
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B23").Value = 'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n'
$ws.Range("C23").Value = 'LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n'
$ws.Range("B24").Value = 'LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n'
$ws.Range("C24").Value = 'LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n'
$ws.Range("B25").Value = 'LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n'
$ws.Range("C25").Value = 'LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n'

# --- Remove the now-obsolete trailing row 26 ---
$ws.Rows.Item(26).Delete()

# --- Reset auto height for label-only rows (12, 17, 22) so no custom height remains ---
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()
